# Implement new reliability-based capacity construction looping through
# hours in constructed timeslice (#232)
#
# Content changes:
#  - About!A1: "MCF Maximum Capacity Factor"
#       -> "BDTPTUMCF Boolean Does This Plant Type Use Maximum Capacity Factor"
#  - BDTPTUMCF!B10 (biomass row): 1 -> 0
#  - Selection/active-cell bookkeeping: BDTPTUMCF!B11 becomes the
#    remembered selection on that sheet, while the About sheet (still the
#    active tab) ends up selected at A1.

$wb = $excel.ActiveWorkbook

$wsAbout = $wb.Worksheets.Item("About")
$wsData  = $wb.Worksheets.Item("BDTPTUMCF")

# Update the descriptive title text on the About sheet.
$wsAbout.Range("A1").Value = "BDTPTUMCF Boolean Does This Plant Type Use Maximum Capacity Factor"

# Turn off "uses maximum capacity factor" for biomass.
$wsData.Range("B10").Value = 0

# Leave the BDTPTUMCF sheet's remembered selection on B11 ...
$wsData.Activate()
$wsData.Range("B11").Select()

# ... then return focus to the About sheet (the tab that was selected
# originally) with its selection reset to A1.
$wsAbout.Activate()
$wsAbout.Range("A1").Select()
